# Added serial run capture data
# ---------------------------------------------------------------------------
# This script reproduces, via Excel COM automation, the edit described by the
# commit "Added serial run capture data" on the "Aggregate Tests" worksheet:
#   - a new column is inserted (H) holding a "N/A" marker for the existing
#     distributed-run rows, with a new header "Serial coWPAtty run time
#     (milliseconds)"
#   - five new rows (7-11) are appended capturing a serial (non-distributed)
#     coWPAtty run for the same five test-data sets
#   - the "linksys_NotInDictionary_..." test-data cells are flagged in red
#   - the selection/view is left on A11
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aggregate Tests")

# ---------------------------------------------------------------------------
# 1. Insert the new "Serial coWPAtty run time (milliseconds)" column.
#    This pushes the existing "Dist Cow..." column's data (col G) one column
#    to the right only from H onward - i.e. old H ("Result") becomes I, while
#    the new, now-empty H is used for the serial-run marker column.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new column H - header + "N/A" for each pre-existing
#    distributed-run row (these runs have no serial-run counterpart value).
#    Write "Serial"/"N/A"/"N/A - cmdline..." text for the first time via the
#    new rows below *before* the header text, so new shared-string entries
#    come out in the same order as the target workbook (Serial, N/A,
#    N/A - cmdline..., Serial coWPAtty run time...).
# ---------------------------------------------------------------------------

# New rows 7-11: the serial (non-distributed) coWPAtty run, one row per test
# data set, mirroring rows 2-6. Columns A/C/I reuse existing shared strings;
# B introduces "Serial" as the first brand-new string.
$ws.Range("A7").Value = "rbeede"
$ws.Range("B7").Value = "Serial"
$ws.Range("C7").Value = "linksys_FirstDictionary_!8zj39le"
$ws.Range("I7").Value = "Correct - Password Found"

$ws.Range("A8").Value = "rbeede"
$ws.Range("B8").Value = "Serial"
$ws.Range("C8").Value = "linksys_MiddleDictionary_korrelie"
$ws.Range("I8").Value = "Correct - Password Found"

$ws.Range("A9").Value = "rbeede"
$ws.Range("B9").Value = "Serial"
$ws.Range("C9").Value = "linksys_LastDictionary_}ttringe"
$ws.Range("I9").Value = "Correct - Password Found"

$ws.Range("A10").Value = "rbeede"
$ws.Range("B10").Value = "Serial"
$ws.Range("C10").Value = "linksys_NotInDictionary_UnknownPassword5763"
$ws.Range("I10").Value = "Correct - No Solution"

$ws.Range("A11").Value = "rbeede"
$ws.Range("B11").Value = "Serial"
$ws.Range("C11").Value = "wireless_Test_invalid_capture"
$ws.Range("I11").Value = "Correct - No Solution"

# Existing distributed-run rows (2-6): new "N/A" marker in column H. This is
# the second brand-new string, and must be introduced here (before the
# "N/A - cmdline..." text below) so shared-string order matches the target.
$ws.Range("H2").Value = "N/A"
$ws.Range("H3").Value = "N/A"
$ws.Range("H4").Value = "N/A"
$ws.Range("H5").Value = "N/A"
$ws.Range("H6").Value = "N/A"

# New rows 7-11, remaining columns: D introduces the third brand-new string
# ("N/A - cmdline..."); E/F/G reuse the "N/A" string already introduced above.
$ws.Range("D7").Value = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E7").Value = "N/A"
$ws.Range("F7").Value = "N/A"
$ws.Range("G7").Value = "N/A"
$ws.Range("H7").Value = 2

$ws.Range("D8").Value = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E8").Value = "N/A"
$ws.Range("F8").Value = "N/A"
$ws.Range("G8").Value = "N/A"
$ws.Range("H8").Value = 2981

$ws.Range("D9").Value = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E9").Value = "N/A"
$ws.Range("F9").Value = "N/A"
$ws.Range("G9").Value = "N/A"
$ws.Range("H9").Value = 5950

$ws.Range("D10").Value = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E10").Value = "N/A"
$ws.Range("F10").Value = "N/A"
$ws.Range("G10").Value = "N/A"

$ws.Range("D11").Value = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E11").Value = "N/A"
$ws.Range("F11").Value = "N/A"
$ws.Range("G11").Value = "N/A"
$ws.Range("H11").Value = 3

# Header for the new column (written last so the "Serial coWPAtty..." string
# is appended after the other new strings above).
$ws.Range("H1").Value = "Serial coWPAtty run time (milliseconds)"

# ---------------------------------------------------------------------------
# 3. Highlight the "not in dictionary" test-data cells in red, in both the
#    distributed run (C5) and the new serial run (C10).
# ---------------------------------------------------------------------------
$ws.Range("C5").Font.Color = 255
$ws.Range("C10").Font.Color = 255

# ---------------------------------------------------------------------------
# 4. Column widths: the new column G keeps the old "Dist Cow..." width, the
#    new column H is sized for the new header, column I keeps the old
#    "Result" width automatically (carried by the Insert above).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 60.59
$ws.Columns.Item(8).ColumnWidth = 41.09

# ---------------------------------------------------------------------------
# 5. Final view state: no frozen/scrolled top-left cell, selection on A11.
# ---------------------------------------------------------------------------
$ws.Range("A11").Select()
